$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 17-18 (Save Game Manager / State Manager), pushing
# Drop Manager .. Record down by two rows, and leaving row 28 as a fresh
# empty row before the trailing DeadLine note row (which lands on row 30).
$ws.Rows("17:18").Insert()

# --- Row 8/9 tweaks ---
$ws.Range("E8").Value = "Being Attack, Attack, setter, getter"
$ws.Range("C9").Value = "Nam 15/02"
$ws.Range("D9").Value = "Hp,Speed,atk"

# --- New row 17: Save Game Manager ---
$ws.Range("A17").Value = "Save Game Manager "
$ws.Range("C17").Value = "Dũng 15/02"
$ws.Range("E17").Value = "Tắt game , thoát từ menu thì save lại lv nhân vật, "
$ws.Range("J17").Value = "wave mấy, hp hiện tại quái, nhân vật, vị trí quái"

# --- New row 18: State Manager ---
$ws.Range("A18").Value = "State Manager"
$ws.Range("D18").Value = "Thêm vào mọi loại quái để checkAvailableStage và thực thi"

# --- Start Scene deadline update (now row 22 after the insert above) ---
$ws.Range("C22").Value = "Nam:15/02/2023"

# --- New row 28: skill-select screen task ---
$ws.Range("A28").Value = "làm màn chọn skill"
$ws.Range("A28").Interior.Color = 16777215
$ws.Range("C28").Value = "Tuấn 15/02"
$ws.Range("E28").Value = "lúc chọn màn chơi thì cho ng chơi chọn 3 kĩ năng "
